# Generate Report for Handoff
# Adds a new file entry (e013cfbb-0a7b-4e95-9277-c9ce08160572.md) as row 9
# to the Overview, zh-cn and de-de tables/sheets.

$wb = $excel.ActiveWorkbook

$fileName   = "e013cfbb-0a7b-4e95-9277-c9ce08160572.md"
$pathName   = "e2e\e013cfbb-0a7b-4e95-9277-c9ce08160572.md"
$commitSha  = "0000000000000000000000000000000000000000"
$baseUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/$fileName"

$zhXlf      = "e013cfbb-0a7b-4e95-9277-c9ce08160572.1c835c3700a86dd172b2a672485b7c1b5a974756.zh-cn.xlf"
$deXlf      = "e013cfbb-0a7b-4e95-9277-c9ce08160572.1c835c3700a86dd172b2a672485b7c1b5a974756.de-de.xlf"

# Empty-cell marker: a lone leading apostrophe forces a text entry whose
# content is the empty string, matching the workbook's existing convention
# of storing blank cells as empty shared strings rather than omitting them.
$EMPTY = "'"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item("Overview")
$rowOverview = $loOverview.ListRows.Add()

$wsOverview.Range("A9").Value = $fileName
$wsOverview.Range("B9").Value = $pathName
$wsOverview.Range("C9").Value = ".md"
$wsOverview.Range("D9").Value = $EMPTY
$wsOverview.Range("D9").Style = "Normal"
$wsOverview.Range("E9").Value = "Ready for handoff"
$wsOverview.Range("F9").Value = "Ready for handoff"
$wsOverview.Range("G9").Value = "2016-09-04 16:48:36"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B9"), $baseUrl, "", "", $pathName)

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item("zh-cn")
$rowZhCn = $loZhCn.ListRows.Add()

$wsZhCn.Range("A9").Value = $fileName
$wsZhCn.Range("B9").Value = ".md"
$wsZhCn.Range("C9").Value = "Ready for handoff"
$wsZhCn.Range("D9").Value = "e2e"
$wsZhCn.Range("E9").Value = "ht"
$wsZhCn.Range("F9").Value = "'False"
$wsZhCn.Range("G9").Value = $zhXlf
$wsZhCn.Range("H9").Value = "2016-09-04 16:48:32"
$wsZhCn.Range("I9").Value = $EMPTY
$wsZhCn.Range("J9").Value = $EMPTY
$wsZhCn.Range("K9").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("L9").Value = $EMPTY
$wsZhCn.Range("M9").Value = "'True"
$wsZhCn.Range("N9").Value = $EMPTY
$wsZhCn.Range("O9").Value = "'False"
$wsZhCn.Range("P9").Value = $EMPTY
$wsZhCn.Range("F9").Style = "Normal"
$wsZhCn.Range("I9:P9").Style = "Normal"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A9"), $baseUrl, "", "", $fileName)

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item("de-de")
$rowDeDe = $loDeDe.ListRows.Add()

$wsDeDe.Range("A9").Value = $fileName
$wsDeDe.Range("B9").Value = ".md"
$wsDeDe.Range("C9").Value = "Ready for handoff"
$wsDeDe.Range("D9").Value = "e2e"
$wsDeDe.Range("E9").Value = "ht"
$wsDeDe.Range("F9").Value = "'False"
$wsDeDe.Range("G9").Value = $deXlf
$wsDeDe.Range("H9").Value = "2016-09-04 16:48:36"
$wsDeDe.Range("I9").Value = $EMPTY
$wsDeDe.Range("J9").Value = $EMPTY
$wsDeDe.Range("K9").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("L9").Value = $EMPTY
$wsDeDe.Range("M9").Value = "'True"
$wsDeDe.Range("N9").Value = $EMPTY
$wsDeDe.Range("O9").Value = "'False"
$wsDeDe.Range("P9").Value = $EMPTY
$wsDeDe.Range("F9").Style = "Normal"
$wsDeDe.Range("I9:P9").Style = "Normal"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A9"), $baseUrl, "", "", $fileName)
